# Update column F (dSF) values for specific rows, per repull/recalculation of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = -6
    10 = -3
    22 = 0
    23 = -6
    24 = -3
    25 = -2
    26 = 0
    28 = -1
    32 = -5
    33 = -5
    34 = -2
    38 = 2
    39 = -3
    42 = 4
    45 = -5
    46 = 3
    49 = -7
    50 = -2
    51 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
